$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy A:F format from row92 to row102 (top-of-group row)
$ws.Range("A92:F92").Copy() | Out-Null
$ws.Range("A102:F102").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Copy A:F format from rows93:F101 to rows103:F109 (middle rows incl bottom row 101->109)
$ws.Range("A93:F101").Copy() | Out-Null
$ws.Range("A103:F111").PasteSpecial(-4122) | Out-Null

Write-Output "done"
